$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rA1 = $wsHoja1.Range("A1")
$text = $rA1.Value2
$text = $text.Replace("1000 Bs = 12.22 = 49024.45 pesos", "1000 Bs = 12.27 = 49534.36 pesos")
$text = $text.Replace("49024.45 pesos = 12.17 = 957.19 Bs", "49534.36 pesos = 12.22 = 976.16 Bs")
$rA1.Value = $text

# --- Update "tasas" sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 81.5
$wsTasas.Range("O10").Value = 4037.05
$wsTasas.Range("N12").Value = 4055
$wsTasas.Range("O12").Value = 79.911
